$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the four timestamp cells (E2:H2) to the new date value.
$newDateSerial = 45840.37710648148
$ws.Range("E2").Value = $newDateSerial
$ws.Range("F2").Value = $newDateSerial
$ws.Range("G2").Value = $newDateSerial
$ws.Range("H2").Value = $newDateSerial

# Update the attachment/file-path text values, replacing the old
# timestamp-based filenames with the new ones.
$oldTag1 = "20250701_130227"
$newTag1 = "20250702_090301"
$oldTag2 = "20250701_130228"
$newTag2 = "20250702_090302"

$k2 = $ws.Range("K2").Value()
$k2 = $k2 -replace $oldTag1, $newTag1 -replace $oldTag2, $newTag2
$ws.Range("K2").Value = $k2

$l2 = $ws.Range("L2").Value()
$l2 = $l2 -replace $oldTag1, $newTag1 -replace $oldTag2, $newTag2
$ws.Range("L2").Value = $l2

$n2 = $ws.Range("N2").Value()
$n2 = $n2 -replace $oldTag1, $newTag1 -replace $oldTag2, $newTag2
$ws.Range("N2").Value = $n2
